# Update countries & provincias Spain
# Applies the daily-data refresh captured in the commit diff:
#  - re-sequences a couple of country-name pairs (their row's stats follow
#    the row, not the label, so the label swap is done alongside the new
#    numbers for that row)
#  - refreshes several countries' case numbers
#  - bumps the "Datos actualizados..." timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp footer ------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 19 de Junio de 2020 a las 05:28"

# --- Country label swap: Surinam / Birmania (rows 160-161) -----------
$ws.Range("A160").Value = "Birmania"
$ws.Range("A161").Value = "Surinam"

# --- Country label swap block: Fiyi/Dominica, Malvinas/Groenlandia,
#     Santa Sede/Islas Turcas y Caicos (rows 202-209) -----------------
$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"
$ws.Range("A206").Value = "Groenlandia"
$ws.Range("A207").Value = "Islas Malvinas"
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("A209").Value = "Santa Sede"

# --- Updated statistics ------------------------------------------------
# Bolivia (row 49)
$ws.Range("B49").Value = 21499
$ws.Range("C49").Value = 814
$ws.Range("E49").Value = 16800
$ws.Range("G49").Value = 18
$ws.Range("H49").Value = 697

# Australia (row 73)
$ws.Range("B73").Value = 7409
$ws.Range("C73").Value = 18
$ws.Range("D73").Value = 6878
$ws.Range("E73").Value = 429

# Sri Lanka (row 104)
$ws.Range("B104").Value = 1947
$ws.Range("C104").Value = 1
$ws.Range("E104").Value = 515

# San Marino (row 139)
$ws.Range("E139").Value = 45
$ws.Range("H139").Value = 42

# Birmania (row 160, label updated above)
$ws.Range("B160").Value = 286
$ws.Range("C160").Value = 23
$ws.Range("D160").Value = 187
$ws.Range("E160").Value = 93
$ws.Range("H160").Value = 6

# Surinam (row 161, label updated above)
$ws.Range("B161").Value = 277
$ws.Range("D161").Value = 74
$ws.Range("E161").Value = 196
$ws.Range("H161").Value = 7

# Islas Turcas y Caicos (row 208, label updated above)
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1

# Santa Sede (row 209, label updated above)
$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0
